$d = $word.ActiveDocument

$replacements = @(
    @("536÷7=", "282÷8="),
    @("922÷4=", "495÷9="),
    @("428÷3=", "567÷9="),
    @("579÷3=", "850÷6="),
    @("653÷9=", "995÷5="),
    @("102÷9=", "842÷9="),
    @("158÷4=", "370÷9="),
    @("129÷2=", "894÷8="),
    @("775÷7=", "692÷2="),
    @("612÷6=", "296÷7="),
    @("889÷6=", "253÷4="),
    @("189÷2=", "823÷2="),
    @("560÷6=", "852÷3="),
    @("436÷2=", "893÷8="),
    @("559÷6=", "264÷5="),
    @("705÷4=", "927÷8="),
    @("772÷4=", "562÷6="),
    @("646÷4=", "646÷7="),
    @("910÷4=", "142÷5="),
    @("696÷8=", "756÷9="),
    @("583÷9=", "716÷4="),
    @("266÷6=", "512÷9="),
    @("370÷3=", "228÷6="),
    @("928÷4=", "352÷9="),
    @("419÷3=", "269÷8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
